$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.126.70'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.478.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.56'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.71%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.478.73'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.07%  '
$ws.Range('E10').Value = '  +2.46%  '
$ws.Range('E11').Value = '  +0.96%  '
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.931.72'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.48'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.073.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  +0.37%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.443.52'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.60'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('E20').Value = '  -3.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.74'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.89%  '
$ws.Range('E22').Value = '  -0.52%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.606.78'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.32%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0909'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '504.59'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.76'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.24'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.54%  '
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.16'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.328'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('E45').Value = '  +1.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '143.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₆0265'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.97%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('E51').Value = '  -1.12%  '
